$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01675466666666667
$ws.Range("H2").Value = 0.050264
$ws.Range("I2").Value = 0.0001854906931657378
$ws.Range("J2").Value = 0.0001854906931657378
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 1.291609266942222
$ws.Range("R2").Value = 11.62448340248
$ws.Range("S2").Value = 0.00004458856210529468
$ws.Range("T2").Value = 0.00004458856210529467
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01675466666666667
$ws.Range("H3").Value = 0.050264
$ws.Range("I3").Value = 0.0001854906931657378
$ws.Range("J3").Value = 0.0001854906931657378
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 1.701939665507556
$ws.Range("R3").Value = 15.317456989568
$ws.Range("S3").Value = 0.00005875386962390287
$ws.Range("T3").Value = 0.00005875386962390286
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01675466666666667
$ws.Range("H4").Value = 0.050264
$ws.Range("I4").Value = 0.0001854906931657378
$ws.Range("J4").Value = 0.0001854906931657378
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 2.379611513016889
$ws.Range("R4").Value = 21.416503617152
$ws.Range("S4").Value = 0.00008214826143654025
$ws.Range("T4").Value = 0.00008214826143654024
$ws.Range("I5").Value = 0.9933938536206305
$ws.Range("J5").Value = 0.9933938536206304
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 6917.202610879296
$ws.Range("R5").Value = 62254.82349791367
$ws.Range("S5").Value = 0.2387936708911015
$ws.Range("T5").Value = 0.2387936708911015
$ws.Range("I6").Value = 0.9933938536206305
$ws.Range("J6").Value = 0.9933938536206304
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.3146558566615664
$ws.Range("T6").Value = 0.3146558566615663
$ws.Range("I7").Value = 0.9933938536206305
$ws.Range("J7").Value = 0.9933938536206304
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.4399443260679626
$ws.Range("T7").Value = 0.4399443260679626
$ws.Range("G8").Value = 0.5799533333333334
$ws.Range("I8").Value = 0.006420655686203657
$ws.Range("J8").Value = 0.006420655686203655
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 44.70832602224445
$ws.Range("R8").Value = 402.3749342002001
$ws.Range("S8").Value = 0.001543407919475529
$ws.Range("T8").Value = 0.001543407919475529
$ws.Range("G9").Value = 0.5799533333333334
$ws.Range("I9").Value = 0.006420655686203657
$ws.Range("J9").Value = 0.006420655686203655
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("Q9").Value = 58.91168125159113
$ws.Range("R9").Value = 530.2051312643201
$ws.Range("S9").Value = 0.002033732046869403
$ws.Range("T9").Value = 0.002033732046869402
$ws.Range("G10").Value = 0.5799533333333334
$ws.Range("I10").Value = 0.006420655686203657
$ws.Range("J10").Value = 0.006420655686203655
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("Q10").Value = 82.3689098964978
$ws.Range("R10").Value = 741.3201890684801
$ws.Range("S10").Value = 0.002843515719858725
$ws.Range("T10").Value = 0.002843515719858724
